$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated projection results values (auto-increase enrollment flag fix)
$ws.Range("C2").Value = 9457
$ws.Range("D2").Value = 8383
$ws.Range("E2").Value = 0.8864333298086073
$ws.Range("F2").Value = 0.8844692973201097
$ws.Range("G2").Value = 0.09581466306536589
$ws.Range("H2").Value = 0.08474512771438723
$ws.Range("I2").Value = 40726470.77878331
$ws.Range("J2").Value = 14166871.48193765
$ws.Range("L2").Value = 14166871.48193765
$ws.Range("M2").Value = 54893342.26072096
$ws.Range("N2").Value = 800181092.3172001
$ws.Range("O2").Value = 782481285.3132
$ws.Range("P2").Value = 0.0177045816477775
$ws.Range("Q2").Value = 0.01810506110221811
$ws.Range("C3").Value = 9643
$ws.Range("D3").Value = 8567
$ws.Range("E3").Value = 0.8884164679041792
$ws.Range("F3").Value = 0.8863025036209394
$ws.Range("G3").Value = 0.09424911249622747
$ws.Range("H3").Value = 0.08353322436945795
$ws.Range("I3").Value = 42515722.80644882
$ws.Range("J3").Value = 14791534.53323031
$ws.Range("L3").Value = 14791534.53323031
$ws.Range("M3").Value = 57307257.33967912
$ws.Range("N3").Value = 837382254.356528
$ws.Range("O3").Value = 819902078.332458
$ws.Range("P3").Value = 0.01766401718722426
$ws.Range("Q3").Value = 0.01804061109750299
$ws.Range("C4").Value = 9835
$ws.Range("D4").Value = 8715
$ws.Range("E4").Value = 0.8861209964412812
$ws.Range("F4").Value = 0.8840535605599513
$ws.Range("G4").Value = 0.09311208234732114
$ws.Range("H4").Value = 0.08231606793030065
$ws.Range("I4").Value = 44356356.26072727
$ws.Range("J4").Value = 15391881.79541372
$ws.Range("L4").Value = 15391881.79541372
$ws.Range("M4").Value = 59748238.05614099
$ws.Range("N4").Value = 875322254.7530119
$ws.Range("O4").Value = 857873306.7470582
$ws.Range("P4").Value = 0.01758424592981109
$ws.Range("Q4").Value = 0.01794190549392158
$ws.Range("C5").Value = 10024
$ws.Range("D5").Value = 8891
$ws.Range("E5").Value = 0.8869712689545092
$ws.Range("F5").Value = 0.8843246469067038
$ws.Range("G5").Value = 0.09183918725606031
$ws.Range("H5").Value = 0.08121565684241421
$ws.Range("I5").Value = 46281693.98751035
$ws.Range("J5").Value = 16021428.19629553
$ws.Range("L5").Value = 16021428.19629553
$ws.Range("M5").Value = 62303122.18380587
$ws.Range("N5").Value = 913377479.2056578
$ws.Range("O5").Value = 895891373.7426846
$ws.Range("P5").Value = 0.01754086186822668
$ws.Range("Q5").Value = 0.01788322632169596
$ws.Range("C6").Value = 10228
$ws.Range("D6").Value = 9091
$ws.Range("E6").Value = 0.8888345717637857
$ws.Range("F6").Value = 0.8865808464989272
$ws.Range("G6").Value = 0.0905776771966571
$ws.Range("H6").Value = 0.08030443372291884
$ws.Range("I6").Value = 48451138.13060883
$ws.Range("J6").Value = 16735055.39106738
$ws.Range("L6").Value = 16735055.39106738
$ws.Range("M6").Value = 65186193.5216762
$ws.Range("N6").Value = 954116855.7025089
$ws.Range("O6").Value = 937584093.1632864
$ws.Range("P6").Value = 0.01753983832383455
$ws.Range("Q6").Value = 0.01786930926354118
